# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for rows 2-25 to their newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 6
    4  = 3
    5  = 3
    6  = 3
    7  = 1
    8  = 4
    9  = 2
    10 = 0
    11 = 6
    12 = 6
    13 = 4
    14 = 6
    15 = 4
    16 = 1
    17 = 3
    18 = 5
    19 = 4
    20 = 5
    21 = 5
    22 = 0
    23 = 1
    24 = 2
    25 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
